$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.000.20"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.650.73"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'309.92"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.3902"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("D8").Value = "'0.3810"
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("D9").Value = "'52.09"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'1.347"
$ws.Range("E10").Value = "  -4.60%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "'0.08451"
$ws.Range("E12").Value = "  -1.81%  "
$ws.Range("D13").Value = "'23.84"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").Value = "'7.058"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").Value = "'8.007"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("E16").Value = "  -2.55%  "
$ws.Range("D17").Value = "1.649.93"
$ws.Range("E17").Value = "  -1.52%  "
$ws.Range("D18").Value = "'94.38"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'0.07008"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'19.68"
$ws.Range("E20").Value = "  -4.52%  "
$ws.Range("D21").Value = "'6.975"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'13.79"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "23.994.22"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "'2.441"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'2.944"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").Value = "'22.07"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "'152.97"
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").Value = "'5.417"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'138.03"
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "'7.933"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").Value = "'2.517"
$ws.Range("E32").Value = "  -0.81%  "
$ws.Range("D33").Value = "1.832.31"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").Value = "'1.017"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("D35").Value = "'0.08049"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").Value = "'6.734"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").Value = "'0.02920"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'10.76"
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "'0.2675"
$ws.Range("E39").Value = "  -3.39%  "
$ws.Range("D40").Value = "'0.09086"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "'0.7600"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").Value = "'13.39"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").Value = "'1.421"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "'16.22"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "'0.6965"
$ws.Range("E45").Value = "  -2.42%  "
$ws.Range("D46").Value = "'2.452"
$ws.Range("E46").Value = "  -3.60%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D49").Value = "'0.08325"
$ws.Range("E49").Value = "  -1.44%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").Value = "'1.229"
$ws.Range("E51").Value = "  -3.57%  "
